# Auto-generated Excel COM-interop script to apply diff changes
# Updates coin rankings table (rows 2-24) and price column (rows 40-51)
# per the commit "Updated symbol list on Sun Dec 25 16:40:08 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D holds price figures that look numeric but are stored as plain
# text. Force Text format on each D cell before writing so Excel keeps the
# literal string (e.g. "243.59") instead of auto-converting to a number.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '243.59'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '23.02'
$ws.Range('B4').Value = 'HuobiToken'
$ws.Range('C4').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.399'
$ws.Range('E4').Value = '3HuobiTokenHT'
$ws.Range('B5').Value = 'Cronos'
$ws.Range('C5').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05934'
$ws.Range('E5').Value = '4CronosCRO'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.398'
$ws.Range('E6').Value = '5GateTokenGT'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8086'
$ws.Range('E7').Value = '6MXTokenMX'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9131'
$ws.Range('E8').Value = '7FTXTokenFTT'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1414'
$ws.Range('E9').Value = '8WazirXWRX'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07423'
$ws.Range('E10').Value = '9MandalaExchangeTokenMDX'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.03316'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03065'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09331'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').Value = 'MCDex'
$ws.Range('C14').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.954'
$ws.Range('E14').Value = '13MCDexMCB'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001573'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04791'
$ws.Range('E16').Value = '15CoinExTokenCET'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.005486'
$ws.Range('E17').Value = '16TigerCashTCH'
$ws.Range('B18').Value = 'HotbitToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.004424'
$ws.Range('E18').Value = '17HotbitTokenHTB'
$ws.Range('B19').Value = 'BitKan'
$ws.Range('C19').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0009861'
$ws.Range('E19').Value = '18BitKanKAN'
$ws.Range('B20').Value = 'NitroEx'
$ws.Range('C20').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.00007804'
$ws.Range('E20').Value = '19NitroExNTX'
$ws.Range('B21').Value = 'LEO'
$ws.Range('C21').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.633'
$ws.Range('E21').Value = '20LEOLEO'
$ws.Range('B22').Value = 'KuCoinToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.431'
$ws.Range('E22').Value = '21KuCoinTokenKCS'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.151'
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.01114'
$ws.Range('E24').Value = '23OneONEBestin24h'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03875'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006213'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1066'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002902'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.006604'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005193'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0005802'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.8347'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002263'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002001'
